$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (avoid numeric auto-conversion)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '98.718.36'
$ws.Range('D3').Value = '3.463.35'
$ws.Range('D5').Value = '261.85'
$ws.Range('D6').Value = '672.98'
$ws.Range('D7').Value = '1.58'
$ws.Range('D8').Value = '0.460'
$ws.Range('D9').Value = '1.14'
$ws.Range('D11').Value = '3.458.89'
$ws.Range('D12').Value = '0.220'
$ws.Range('D13').Value = '42.68'
$ws.Range('D14').Value = '0.0000274'
$ws.Range('D15').Value = '6.24'
$ws.Range('D16').Value = '98.767.24'
$ws.Range('D17').Value = '4.118.43'
$ws.Range('D18').Value = '8.33'
$ws.Range('D19').Value = '3.471.59'
$ws.Range('D20').Value = '17.53'
$ws.Range('D21').Value = '3.62'
$ws.Range('D22').Value = '535.97'
$ws.Range('D24').Value = '0.478'
$ws.Range('D25').Value = '0.0000218'
$ws.Range('D26').Value = '6.44'
$ws.Range('D27').Value = '103.37'
$ws.Range('D28').Value = '12.97'
$ws.Range('D29').Value = '3.661.35'
$ws.Range('D30').Value = '0.153'
$ws.Range('D31').Value = '11.58'
$ws.Range('D32').Value = '0.200'
$ws.Range('D33').Value = '0.999'
$ws.Range('D34').Value = '0.576'
$ws.Range('D35').Value = '0.998'
$ws.Range('D36').Value = '30.56'
$ws.Range('D37').Value = '2.24'
$ws.Range('D38').Value = '8.05'
$ws.Range('D39').Value = '0.162'
$ws.Range('D40').Value = '538.56'
$ws.Range('D41').Value = '1.42'
$ws.Range('D42').Value = '24.78'
$ws.Range('D44').Value = '0.863'
$ws.Range('D45').Value = '3.51'
$ws.Range('D46').Value = '3.74'
$ws.Range('D47').Value = '8.33'
$ws.Range('D50').Value = '2.12'
$ws.Range('D51').Value = '5.30'

# Restore default (unstyled) formatting on the price column so styling matches original
$ws.Range('D2:D51').ClearFormats()

# Update Coin / Link / Volume(1h) columns
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('E3').Value = '  +4.73%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  +2.23%  '
$ws.Range('E6').Value = '  +8.05%  '
$ws.Range('E7').Value = '  +9.39%  '
$ws.Range('E8').Value = '  +14.59%  '
$ws.Range('E9').Value = '  +24.76%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  +4.63%  '
$ws.Range('E12').Value = '  +10.53%  '
$ws.Range('E13').Value = '  +9.78%  '
$ws.Range('E14').Value = '  +10.92%  '
$ws.Range('E15').Value = '  +13.93%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('E17').Value = '  +5.05%  '
$ws.Range('E18').Value = '  +34.21%  '
$ws.Range('E19').Value = '  +4.95%  '
$ws.Range('E20').Value = '  +15.27%  '
$ws.Range('E21').Value = '  +3.29%  '
$ws.Range('E22').Value = '  +11.09%  '
$ws.Range('E23').Value = '  +14.22%  '
$ws.Range('B24').Value = 'Stellar'
$ws.Range('C24').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E24').Value = '  +60.60%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('E25').Value = '  +6.88%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E26').Value = '  +15.19%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E27').Value = '  +16.71%  '
$ws.Range('E28').Value = '  +9.34%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('E29').Value = '  +5.32%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E30').Value = '  +16.73%  '
$ws.Range('E31').Value = '  +17.15%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E32').Value = '  +7.59%  '
$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('B34').Value = 'PolygonEcosystemToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('E34').Value = '  +25.94%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E36').Value = '  +10.47%  '
$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E37').Value = '  +15.52%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('E38').Value = '  +12.19%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E39').Value = '  +9.54%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('E40').Value = '  +9.74%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('E41').Value = '  +14.94%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('E43').Value = '  +34.82%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E44').Value = '  +7.60%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E45').Value = '  +11.44%  '
$ws.Range('B46').Value = 'MantraDAO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E47').Value = '  +16.61%  '
$ws.Range('E48').Value = '  +19.11%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('E50').Value = '  +11.58%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E51').Value = '  +14.52%  '
